$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D (Price) and E (Volume) column cells are stored as text in the
# source sheet (e.g. "3.314.39", "  +6.03%  "). Prefix with a leading
# apostrophe so Excel keeps them as text instead of coercing to numbers
# (which would lose formatting / introduce floating point drift).

$ws.Range('D2').Value = "'63.952.68"
$ws.Range('E2').Value = "'  +1.46%  "
$ws.Range('D3').Value = "'3.314.39"
$ws.Range('E3').Value = "'  +6.03%  "
$ws.Range('E4').Value = "'  +0.06%  "
$ws.Range('D5').Value = "'598.78"
$ws.Range('E5').Value = "'  +0.86%  "
$ws.Range('D6').Value = "'143.36"
$ws.Range('E6').Value = "'  +5.13%  "
$ws.Range('E7').Value = "'  -0.05%  "
$ws.Range('D8').Value = "'3.312.84"
$ws.Range('E8').Value = "'  +6.26%  "
$ws.Range('E10').Value = "'  +2.89%  "
$ws.Range('D11').Value = "'5.51"
$ws.Range('E11').Value = "'  +5.72%  "
$ws.Range('D12').Value = "'0.474"
$ws.Range('E12').Value = "'  +3.70%  "
$ws.Range('E13').Value = "'  +1.26%  "
$ws.Range('D14').Value = "'34.81"
$ws.Range('E14').Value = "'  +1.78%  "
$ws.Range('D15').Value = "'3.861.29"
$ws.Range('E15').Value = "'  +6.15%  "
$ws.Range('E16').Value = "'  +1.16%  "
$ws.Range('D17').Value = "'3.311.59"
$ws.Range('E17').Value = "'  +5.79%  "
$ws.Range('D18').Value = "'64.025.09"
$ws.Range('E18').Value = "'  +1.67%  "
$ws.Range('D19').Value = "'6.91"
$ws.Range('E19').Value = "'  +3.07%  "
$ws.Range('D20').Value = "'482.18"
$ws.Range('E20').Value = "'  +1.59%  "
$ws.Range('E21').Value = "'  +0.21%  "
$ws.Range('D22').Value = "'0.741"
$ws.Range('E22').Value = "'  +5.96%  "
$ws.Range('E23').Value = "'  +3.77%  "
$ws.Range('D24').Value = "'13.61"
$ws.Range('E24').Value = "'  +4.55%  "
$ws.Range('D25').Value = "'84.72"
$ws.Range('E25').Value = "'  -2.65%  "
$ws.Range('E26').Value = "'  +0.28%  "
$ws.Range('E27').Value = "'  +2.17%  "
$ws.Range('D28').Value = "'7.28"
$ws.Range('E28').Value = "'  +1.48%  "
$ws.Range('E29').Value = "'  -0.19%  "
$ws.Range('E30').Value = "'  +3.74%  "
$ws.Range('E31').Value = "'  +5.24%  "
$ws.Range('D32').Value = "'29.14"
$ws.Range('E32').Value = "'  +9.00%  "
$ws.Range('E33').Value = "'  -0.25%  "
$ws.Range('E34').Value = "'  +1.15%  "
$ws.Range('E35').Value = "'  +2.43%  "
$ws.Range('E36').Value = "'  +3.10%  "
$ws.Range('E37').Value = "'  +2.40%  "
$ws.Range('D38').Value = "'0.0₃0760"
$ws.Range('E38').Value = "'  +8.10%  "
$ws.Range('D39').Value = "'0.0400"
$ws.Range('E39').Value = "'  +3.44%  "
$ws.Range('D40').Value = "'432.82"
$ws.Range('E40').Value = "'  +2.91%  "
$ws.Range('D41').Value = "'3.043.48"
$ws.Range('E41').Value = "'  +5.53%  "
$ws.Range('D42').Value = "'8.45"
$ws.Range('E42').Value = "'  +2.41%  "
$ws.Range('E43').Value = "'  +3.36%  "
$ws.Range('E44').Value = "'  -1.49%  "
$ws.Range('D45').Value = "'0.268"
$ws.Range('E45').Value = "'  +2.27%  "
$ws.Range('E46').Value = "'  +4.18%  "
$ws.Range('D47').Value = "'26.54"
$ws.Range('E47').Value = "'  +3.07%  "
$ws.Range('E49').Value = "'  +2.47%  "
$ws.Range('D50').Value = "'35.56"
$ws.Range('E50').Value = "'  +11.97%  "
$ws.Range('E51').Value = "'  +1.77%  "
